$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "AYP Met?" to column O, row 1, matching the style of existing headers
$ws.Range("O1").Value = "AYP Met?"
$ws.Range("O1").Font.Bold = $true
$ws.Range("O1").Font.Size = 16

# Update the selection / view to reflect where the user ended up after adding the column
$ws.Range("O2").Select()
